$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---
$ws.Range("H2").NumberFormat = "0.00E+00"
$ws.Range("H2").Value2 = 0.58986973762512196
$ws.Range("I2").HorizontalAlignment = 1
$ws.Range("I2").Value2 = 0.50197410583496005

# --- Row 3 ---
$ws.Range("H3").Value2 = 0.83108949661254805
$ws.Range("I3").HorizontalAlignment = 1
$ws.Range("I3").Value2 = 0.53224802017211903

# --- Row 4 ---
$ws.Range("H4").Value2 = 0.70217561721801702
$ws.Range("I4").HorizontalAlignment = 1
$ws.Range("I4").Value2 = 0.57619571685791005

# --- Row 5 ---
$ws.Range("H5").Value2 = 0.68948197364807096
$ws.Range("I5").NumberFormat = "0.00E+00"

# --- Row 6 ---
$ws.Range("H6").NumberFormat = "0.00E+00"
$ws.Range("H6").Value2 = 14.1812629699707
$ws.Range("I6").NumberFormat = "0.00E+00"

# --- Row 7 ---
$ws.Range("H7").NumberFormat = "0.00E+00"
$ws.Range("H7").Value2 = 21.4208245277404
$ws.Range("I7").HorizontalAlignment = 1
$ws.Range("I7").NumberFormat = "0.00E+00"
$ws.Range("I7").Value2 = 14.988908767700099

# --- Row 8 ---
$ws.Range("H8").NumberFormat = "0.00E+00"
$ws.Range("H8").Value2 = 15.209624767303399
$ws.Range("I8").NumberFormat = "0.00E+00"

# --- Row 9 ---
$ws.Range("H9").NumberFormat = "0.00E+00"
$ws.Range("H9").Value2 = 11.3813378810882
$ws.Range("I9").NumberFormat = "0.00E+00"

# --- Row 10 ---
$ws.Range("H10").NumberFormat = "0.00E+00"
$ws.Range("H10").Value2 = 15.572923898696899
$ws.Range("I10").HorizontalAlignment = 1
$ws.Range("I10").NumberFormat = "0.00E+00"
$ws.Range("I10").Value2 = 2.2686486244201598

# --- Row 11 ---
$ws.Range("H11").NumberFormat = "0.00E+00"
$ws.Range("H11").Value2 = 6.68387603759765
$ws.Range("I11").NumberFormat = "0.00E+00"

# --- Row 12 ---
$ws.Range("H12").NumberFormat = "0.00E+00"
$ws.Range("H12").Value2 = 9.1830027103424001
$ws.Range("I12").NumberFormat = "0.00E+00"

# --- Row 13 ---
$ws.Range("H13").NumberFormat = "0.00E+00"
$ws.Range("H13").Value2 = 13.751554250717099
$ws.Range("I13").HorizontalAlignment = 1
$ws.Range("I13").NumberFormat = "0.00E+00"
$ws.Range("I13").Value2 = 2.10360360145568

# --- Row 14 ---
$ws.Range("H14").NumberFormat = "0.00E+00"
$ws.Range("H14").Value2 = 18.754692316055198
$ws.Range("I14").HorizontalAlignment = 1
$ws.Range("I14").NumberFormat = "0.00E+00"
$ws.Range("I14").Value2 = 1.85164093971252

# --- Row 15 ---
$ws.Range("H15").NumberFormat = "0.00E+00"
$ws.Range("H15").Value2 = 15.477213382720899
$ws.Range("I15").NumberFormat = "0.00E+00"

# --- Row 16 ---
$ws.Range("H16").NumberFormat = "0.00E+00"
$ws.Range("H16").Value2 = 21.049712896347
$ws.Range("I16").HorizontalAlignment = 1
$ws.Range("I16").NumberFormat = "0.00E+00"
$ws.Range("I16").Value2 = 3.2149789333343501

# --- Row 17 ---
$ws.Range("H17").NumberFormat = "0.00E+00"
$ws.Range("H17").Value2 = 4.2755711078643799
$ws.Range("I17").NumberFormat = "0.00E+00"
$ws.Range("I17").HorizontalAlignment = 1
$ws.Range("I17").Value2 = 6.8606407642364502

# --- Selection ---
$ws.Range("G18").Select()
